$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.526.39'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.467.41'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").Value = '''315.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").Value = '''91.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.64%  '
$ws.Range("D7").Value = '''0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").Value = '''0.515'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.33%  '
$ws.Range("D10").Value = '''32.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.35%  '
$ws.Range("D11").Value = '''0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '2.848.49'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '''6.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").Value = '''15.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").Value = '2.465.88'
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '''0.776'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").Value = '41.547.41'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = '''6.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").Value = '0.0₃0943'
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").Value = '''71.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").Value = '''11.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '''236.67'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = '''24.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("E28").Value = '  -0.96%  '
$ws.Range("D29").Value = '''9.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("D30").Value = '''35.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.39%  '
$ws.Range("D31").Value = '''155.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.39%  '
$ws.Range("D32").Value = '''5.44'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").Value = '''2.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '''0.0760'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("D35").Value = '''17.14'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.44%  '
$ws.Range("D36").Value = '''2.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = '''2.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.56%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.103'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = '''0.114'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("D40").Value = '''1.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.30%  '
$ws.Range("D41").Value = '''3.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("D43").Value = '1.943.28'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").Value = '''0.0283'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("D45").Value = '''18.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.76%  '
$ws.Range("D46").Value = '''2.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.94%  '
$ws.Range("D47").Value = '''9.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.38%  '
$ws.Range("D48").Value = '2.708.14'
$ws.Range("D49").Value = '''97.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '''67.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.77%  '
$ws.Range("D51").Value = '''52.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.61%  '
